$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = '@'
$ws.Cells.Item(2, 4).Value = '25.972.02'
$ws.Cells.Item(2, 5).Value = '  -0.59%  '
$ws.Cells.Item(3, 4).NumberFormat = '@'
$ws.Cells.Item(3, 4).Value = '1.643.34'
$ws.Cells.Item(3, 5).Value = '  -1.41%  '
$ws.Cells.Item(4, 5).Value = '  -0.06%  '
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '214.78'
$ws.Cells.Item(5, 5).Value = '  +2.26%  '
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '0.5215'
$ws.Cells.Item(6, 5).Value = '  -0.01%  '
$ws.Cells.Item(7, 4).NumberFormat = '@'
$ws.Cells.Item(7, 4).Value = '1.001'
$ws.Cells.Item(7, 5).Value = '  -0.12%  '
$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '0.2605'
$ws.Cells.Item(8, 5).Value = '  +0.16%  '
$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '0.06350'
$ws.Cells.Item(9, 5).Value = '  +0.42%  '
$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '20.72'
$ws.Cells.Item(10, 5).Value = '  -1.68%  '
$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '0.07672'
$ws.Cells.Item(11, 5).Value = '  +1.89%  '
$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '1.644.08'
$ws.Cells.Item(12, 5).Value = '  -1.06%  '
$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '4.423'
$ws.Cells.Item(13, 5).Value = '  +0.29%  '
$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '1.868.23'
$ws.Cells.Item(14, 5).Value = '  -1.29%  '
$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '0.5522'
$ws.Cells.Item(15, 5).Value = '  +1.92%  '
$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '0.0₅8273'
$ws.Cells.Item(16, 5).Value = '  +3.27%  '
$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '64.69'
$ws.Cells.Item(17, 5).Value = '  -2.50%  '
$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '25.980.60'
$ws.Cells.Item(18, 5).Value = '  -0.79%  '
$ws.Cells.Item(19, 5).Value = '  +0.01%  '
$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '4.702'
$ws.Cells.Item(20, 5).Value = '  -0.61%  '
$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '189.10'
$ws.Cells.Item(21, 5).Value = '  +1.11%  '
$ws.Cells.Item(22, 5).Value = '  -0.83%  '
$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '6.256'
$ws.Cells.Item(23, 5).Value = '  +0.52%  '
$ws.Cells.Item(24, 5).Value = '  -0.04%  '
$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '144.29'
$ws.Cells.Item(25, 5).Value = '  -3.55%  '
$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '0.1237'
$ws.Cells.Item(26, 5).Value = '  -0.06%  '
$ws.Cells.Item(27, 5).Value = '  -0.55%  '
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '15.89'
$ws.Cells.Item(28, 5).Value = '  +1.07%  '
$ws.Cells.Item(29, 5).Value = '  +2.22%  '
$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '0.05919'
$ws.Cells.Item(30, 5).Value = '  -5.79%  '
$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '1.262'
$ws.Cells.Item(31, 5).Value = '  -1.07%  '
$ws.Cells.Item(32, 5).Value = '  -0.35%  '
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '3.393'
$ws.Cells.Item(33, 5).Value = '  -2.95%  '
$ws.Cells.Item(34, 5).Value = '  +0.67%  '
$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '0.9931'
$ws.Cells.Item(35, 5).Value = '  -0.66%  '
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '2.394'
$ws.Cells.Item(36, 5).Value = '  +0.08%  '
$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '2.747'
$ws.Cells.Item(37, 5).Value = '  -0.56%  '
$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '0.5632'
$ws.Cells.Item(38, 5).Value = '  -5.92%  '
$ws.Cells.Item(39, 5).Value = '  -0.32%  '
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '5.868'
$ws.Cells.Item(40, 5).Value = '  -3.04%  '
$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '0.8520'
$ws.Cells.Item(41, 5).Value = '  -1.18%  '
$ws.Cells.Item(42, 5).Value = '  -0.16%  '
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '1.030.95'
$ws.Cells.Item(43, 5).Value = '  -7.10%  '
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '98.80'
$ws.Cells.Item(44, 5).Value = '  -1.84%  '
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '1.791.43'
$ws.Cells.Item(45, 5).Value = '  -1.43%  '
$ws.Cells.Item(46, 2).Value = 'Aave'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '55.62'
$ws.Cells.Item(46, 5).Value = '  +0.61%  '
$ws.Cells.Item(47, 2).Value = 'Frax'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '1.003'
$ws.Cells.Item(47, 5).Value = '  +0.11%  '
$ws.Cells.Item(48, 2).Value = 'EnergySwap'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '8.035'
$ws.Cells.Item(48, 5).Value = '  -0.47%  '
$ws.Cells.Item(49, 2).Value = 'Cronos'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '0.05145'
$ws.Cells.Item(49, 5).Value = '  -1.89%  '
$ws.Cells.Item(50, 2).Value = 'Mantle'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '0.4219'
$ws.Cells.Item(50, 5).Value = '  -0.31%  '
$ws.Cells.Item(51, 2).Value = 'Aptos'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '5.891'
$ws.Cells.Item(51, 5).Value = '  +0.20%  '
